# #5: insurance, claim, debt, investment done
#
# The original "保險" (insurance) and "債務" (debt) sheets only carried a
# handful of raw columns. Bring them in line with the "土地" (land) sheet by
# adding the full common metadata columns (property_category/category/date/
# legislator_name/legislator_id/source_file/index, plus a "total" and
# "owner"/"debtor"/"species" breakout for the debt sheet) and giving every
# row an explicit company/name resp. species/debtor header.

$wb = $excel.ActiveWorkbook
$ws_land = $wb.Worksheets.Item(1)
$ws_ins  = $wb.Worksheets.Item(2)
$ws_debt = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet 2 - 保險 (insurance)
# ---------------------------------------------------------------------
# Header row: relabel existing header cells and extend with the shared
# metadata headers (copy formatting from the land sheet's header row so the
# new cells pick up the bold/bordered header style).
$ws_land.Range("I1:O1").Copy() | Out-Null
$ws_ins.Range("E1:K1").PasteSpecial(-4122) | Out-Null

$ws_ins.Range("B1").Value = "company"
$ws_ins.Range("C1").Value = "name"
$ws_ins.Range("D1").Value = "owner"
$ws_ins.Range("E1").Value = "property_category"
$ws_ins.Range("F1").Value = "category"
$ws_ins.Range("G1").Value = "date"
$ws_ins.Range("H1").Value = "legislator_name"
$ws_ins.Range("I1").Value = "legislator_id"
$ws_ins.Range("J1").Value = "source_file"
$ws_ins.Range("K1").Value = "index"

# Data rows 2-5: copy the data-row cell format into the newly used columns,
# then fill in the values. B/C/D already hold the correct company/name/owner
# values, so only E:K are new.
$insRows = @(2, 3, 4, 5)
foreach ($r in $insRows) {
    $ws_land.Range("I2:O2").Copy() | Out-Null
    $ws_ins.Range("E${r}:K${r}").PasteSpecial(-4122) | Out-Null

    $ws_ins.Range("E$r").Value = "insurance"
    $ws_ins.Range("F$r").Value = "normal"
    $ws_ins.Range("G$r").Value = "2011-12-30"
    $ws_ins.Range("H$r").Value = "陳淑慧"
    $ws_ins.Range("I$r").Value = 1720
    $ws_ins.Range("J$r").Value = "tmp503c1"
}
$ws_ins.Range("K2").Value = 101
$ws_ins.Range("K3").Value = 102
$ws_ins.Range("K4").Value = 103
$ws_ins.Range("K5").Value = 104

# ---------------------------------------------------------------------
# Sheet 3 - 債務 (debt)
# ---------------------------------------------------------------------
# The sheet had been laid out one column short (no "species" column, with a
# duplicated "借款" value standing in for "debtor"). Rebuild the header and
# the two data rows with the correct species/debtor/owner/total/register_*
# breakout, followed by the same shared metadata columns as the other
# sheets.
$ws_land.Range("I1:O1").Copy() | Out-Null
$ws_debt.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws_debt.Range("B1").Value = "species"
$ws_debt.Range("C1").Value = "debtor"
$ws_debt.Range("D1").Value = "owner"
$ws_debt.Range("E1").Value = "total"
$ws_debt.Range("F1").Value = "register_date"
$ws_debt.Range("G1").Value = "register_reason"
$ws_debt.Range("H1").Value = "property_category"
$ws_debt.Range("I1").Value = "category"
$ws_debt.Range("J1").Value = "date"
$ws_debt.Range("K1").Value = "legislator_name"
$ws_debt.Range("L1").Value = "legislator_id"
$ws_debt.Range("M1").Value = "source_file"
$ws_debt.Range("N1").Value = "index"

# Row 2: 元大商銀 loan
$ws_land.Range("I2:O2").Copy() | Out-Null
$ws_debt.Range("H2:N2").PasteSpecial(-4122) | Out-Null

$ws_debt.Range("B2").Value = "借款"
$ws_debt.Range("C2").Value = "林南生"
$ws_debt.Range("D2").Value = "元大商銀臺南市中西區民生路"
$ws_debt.Range("E2").Value = 2367290
$ws_debt.Range("F2").Value = "96年10月15日"
$ws_debt.Range("G2").Value = "借款"
$ws_debt.Range("H2").Value = "debt"
$ws_debt.Range("I2").Value = "normal"
$ws_debt.Range("J2").Value = "2011-12-30"
$ws_debt.Range("K2").Value = "陳淑慧"
$ws_debt.Range("L2").Value = 1720
$ws_debt.Range("M2").Value = "tmp503c1"
$ws_debt.Range("N2").Value = 127

# Row 3: 安泰商業銀行 loan
$ws_land.Range("I2:O2").Copy() | Out-Null
$ws_debt.Range("H3:N3").PasteSpecial(-4122) | Out-Null

$ws_debt.Range("B3").Value = "借款"
$ws_debt.Range("C3").Value = "林南生"
$ws_debt.Range("D3").Value = "安泰商業銀行臺南市中西區中山路"
$ws_debt.Range("E3").Value = 3067394
$ws_debt.Range("F3").Value = "94年08月02日"
$ws_debt.Range("G3").Value = "借款"
$ws_debt.Range("H3").Value = "debt"
$ws_debt.Range("I3").Value = "normal"
$ws_debt.Range("J3").Value = "2011-12-30"
$ws_debt.Range("K3").Value = "陳淑慧"
$ws_debt.Range("L3").Value = 1720
$ws_debt.Range("M3").Value = "tmp503c1"
$ws_debt.Range("N3").Value = 128
